# Generate Report for Handback
# Update generation/handoff/handback timestamps and priority flags that
# resulted from re-running the handback report generation.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-25 04:16:27"
$wsOverview.Range("G4").Value = "2016-08-25 04:16:27"

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-25 04:16:23"
$wsZhCn.Range("H4").Value = "2016-08-25 04:16:23"
$wsZhCn.Range("K3").Value = "2016-08-25 04:16:39"
$wsZhCn.Range("K4").Value = "2016-08-25 04:16:39"

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-25 04:16:27"
$wsDeDe.Range("H4").Value = "2016-08-25 04:16:27"
$wsDeDe.Range("K3").Value = "2016-08-25 04:16:46"
$wsDeDe.Range("K4").Value = "2016-08-25 04:16:46"
